$wb = $excel.ActiveWorkbook

# --- Update selection on Student_preferences sheet ---
$wsStudent = $wb.Worksheets.Item("Student_preferences")
$wsStudent.Range("D2").Select()

# --- Update selection on Supervisor_preferences sheet ---
$wsSupervisor = $wb.Worksheets.Item("Supervisor_preferences")
$wsSupervisor.Range("B1").Select()

# --- Add the new COMMENTS sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsComments = $wb.Worksheets.Add($null, $lastSheet)
$wsComments.Name = "COMMENTS"

$wsComments.Range("A1").Value = "Try this one with different supervisor weights."
$wsComments.Range("A2").Value = "With supervisor weight 0, students S1 and S2 get their preference of project P1."
$wsComments.Range("A3").Value = "But with supervisors weights e.g. 0.25-5, they get pulled towards project P2."
$wsComments.Range("B4").Value = "… because supervisor P2, who cares, is giving them supervisor dissatisfaction scores of 1 and 2..."
$wsComments.Range("B5").Value = "… and supervisor P1, who doesn’t care, is by default giving all students scores of 5.5"
$wsComments.Range("B6").Value = "… so the supervisors are “happier” if both students are allocated to P2."
$wsComments.Range("A7").Value = "(And with weight 0.2, student S2 gets pulled to project P2 but student S1 stays with project P1.)"
$wsComments.Range("A9").Value = "This is all logically consistent."
$wsComments.Range("A10").Value = "It’s not clear there’s a better method:"
$wsComments.Range("A11").Value = "- if supervisors who don't get a preference allocate 0 for each student (not 5.5), that pulls students towards them, rather than away"
$wsComments.Range("A12").Value = "- if we allocated them the mean expressed preference for supervisors who expressed a preference, that cannot be done (or is not consistent) when some supervisors rank 2 students, some rank 4, etc."
$wsComments.Range("A14").Value = "The only question is how much the supervisors should be allowed to influence things (definitely not zero; as a matter of policy less than student preferences)."
$wsComments.Range("A16").Value = "The other question is whether this causes any problems in practice."
$wsComments.Range("A18").Value = "(NB A previous computer-based attempt in PDF failed, i.e. caused unhappiness and needed redoing manually, but we have reason to believe it wasn’t very sophisticated computationally, using ?LiveCode.)"
$wsComments.Range("B19").Value = "https://en.wikipedia.org/wiki/LiveCode"
$wsComments.Range("B20").Value = "https://livecode.org/"
$wsComments.Range("B21").Value = "… which, while not entirely incapable, does not look like it has an integer programming library, and I would not want to create one for scratch for it as it looks dreadful."
$wsComments.Range("B22").Value = "… Google: “livecode” “integer programming” – 0 hits"

$wsComments.Range("A19").Select()
